{"js": "// Resize the two columns of the Code Review form's first table\n// (5613/3993 dxa -> 6220/3386 dxa) and update the developer lab-group\n// cell's wording.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// `TableCell.columnWidth` is column-wide: writing it resizes every cell in\n// that column plus the table's <w:gridCol>, so a single write per column\n// keeps the grid consistent (matches the diff's gridCol + 4x tcW updates).\n// Widths are expressed in points; dxa / 20 = points.\ntable.getCell(0, 0).columnWidth = 6220 / 20; // 311.0 pt  (was 5613 dxa / 280.65 pt)\ntable.getCell(0, 1).columnWidth = 3386 / 20; // 169.3 pt  (was 3993 dxa / 199.65 pt)\n\n// Update the \"Developer's Lab Assignment Group\" row label.\nconst oldText = \"Developer\\u2019s Lab Assignment Group (A, B, or C)\";\nconst newText = \"Developer\\u2019s lab version (A, B, or C) and lab partner group name\";\nconst matches = body.search(oldText, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  matches.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Resize the two columns of the Code Review form's first table\n# (5613/3993 dxa -> 6220/3386 dxa, i.e. 280.65/199.65 pt -> 311.0/169.3 pt)\n# and update the developer lab-group cell's wording.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# Table.Columns(i).Width resizes every cell in that column plus the\n# table's grid column, keeping the table's <w:tblGrid>/<w:tcW> consistent.\n$t.Columns(1).Width = 311.0\n$t.Columns(2).Width = 169.3\n\n# Update the \"Developer's Lab Assignment Group\" row label.\n$rsquo = [char]0x2019\n$findText = \"Developer\" + $rsquo + \"s Lab Assignment Group (A, B, or C)\"\n$replaceText = \"Developer\" + $rsquo + \"s lab version (A, B, or C) and lab partner group name\"\n\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Replacement.ClearFormatting()\n$findRange.Find.Text = $findText\n$findRange.Find.Replacement.Text = $replaceText\n$findRange.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n"}
